$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.418.95'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '2.595.42'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '''535.17'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '''142.20'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''0.584'
$ws.Range("E8").Value = '  +3.08%  '
$ws.Range("D9").Value = '2.601.92'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '''6.81'
$ws.Range("E10").Value = '  +3.34%  '
$ws.Range("D11").Value = '''0.0998'
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").Value = '''0.332'
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").Value = '3.071.31'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '58.385.92'
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").Value = '''20.72'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '2.595.43'
$ws.Range("E17").Value = '  -3.38%  '
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '''4.40'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '''333.70'
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("D21").Value = '''10.09'
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = '''6.18'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '''66.96'
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").Value = '''0.420'
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  -3.15%  '
$ws.Range("D28").Value = '''7.07'
$ws.Range("E28").Value = '  -2.20%  '
$ws.Range("B29").Value = 'USDe'
$ws.Range("C29").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0731'
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("D31").Value = '''1.64'
$ws.Range("E31").Value = '  -1.62%  '
$ws.Range("D32").Value = '''5.88'
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("D33").Value = '''153.57'
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").Value = '''18.80'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '''3.89'
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("D36").Value = '''37.02'
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("D37").Value = '''1.10'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").Value = '''0.832'
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = '''0.820'
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("D40").Value = '''1.41'
$ws.Range("E40").Value = '  -3.07%  '
$ws.Range("D41").Value = '''3.59'
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("D42").Value = '''281.67'
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value = '''10.69'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.590'
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("D46").Value = '''0.0949'
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = '''19.00'
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("D48").Value = '''0.0529'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = '''0.0226'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '1.932.04'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").Value = '''4.45'
$ws.Range("E51").Value = '  -1.56%  '
